# Applies the stock-count/value corrections described in the commit diff:
# quantities (col F) and computed values (col G) shrink for a set of items
# (stock recount), a few rows swap their B/D/E/F/G contents (re-sorted
# pairs with identical item names), and the "Sub Total:" / "Grand Total:"
# rows in column B are updated to the new summed totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 58
$ws.Range("G6").Value = 1733.04
$ws.Range("B10").Value = 26947.48
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 1105.92
$ws.Range("F31").Value = 29
$ws.Range("G31").Value = 1024.28
$ws.Range("B32").Value = 12326.83
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 404.85
$ws.Range("F86").Value = 53
$ws.Range("G86").Value = 6649.91
$ws.Range("B90").Value = 170237.58
$ws.Range("F115").Value = 188
$ws.Range("G115").Value = 18200.28
$ws.Range("B117").Value = 12041.34
$ws.Range("F144").Value = 974
$ws.Range("G144").Value = 8230.299999999999
$ws.Range("F146").Value = 17
$ws.Range("G146").Value = 1431.23
$ws.Range("B147").Value = 12841.55
$ws.Range("F150").Value = 31
$ws.Range("G150").Value = 1441.19
$ws.Range("F153").Value = 18
$ws.Range("G153").Value = 836.8200000000001
$ws.Range("B156").Value = 30348.13
$ws.Range("F160").Value = 13
$ws.Range("G160").Value = 1246.83
$ws.Range("B161").Value = 1286.58
$ws.Range("F164").Value = 61
$ws.Range("G164").Value = 6992.43
$ws.Range("B175").Value = 26737.72
$ws.Range("F197").Value = 18
$ws.Range("G197").Value = 1117.08
$ws.Range("F207").Value = 15
$ws.Range("G207").Value = 403.35
$ws.Range("F214").Value = 42
$ws.Range("G214").Value = 3683.4
$ws.Range("B216").Value = 35899.13
$ws.Range("F218").Value = 4
$ws.Range("G218").Value = 864.88
$ws.Range("F233").Value = 114
$ws.Range("G233").Value = 5430.96
$ws.Range("F234").Value = 37
$ws.Range("G234").Value = 1898.84
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("F250").Value = 6
$ws.Range("G250").Value = 2966.28
$ws.Range("F255").Value = 523
$ws.Range("G255").Value = 89605.59
$ws.Range("F256").Value = 266
$ws.Range("G256").Value = 40211.22
$ws.Range("B260").Value = 173220.47
$ws.Range("F280").Value = 130
$ws.Range("G280").Value = 21988.2
$ws.Range("F293").Value = 30
$ws.Range("G293").Value = 2109.6
$ws.Range("F294").Value = 26
$ws.Range("G294").Value = 1855.36
$ws.Range("F302").Value = 38
$ws.Range("G302").Value = 8013.82
$ws.Range("F303").Value = 23
$ws.Range("G303").Value = 4850.47
$ws.Range("B304").Value = 167080
$ws.Range("B322").Value = 47097
$ws.Range("D322").Value = 112.28
$ws.Range("E322").Value = 134.16
$ws.Range("F322").Value = 15
$ws.Range("G322").Value = 1684.2
$ws.Range("B323").Value = 58047
$ws.Range("D323").Value = 105.54
$ws.Range("E323").Value = 126.1
$ws.Range("F323").Value = 39
$ws.Range("G323").Value = 4116.06
$ws.Range("F334").Value = 190
$ws.Range("G334").Value = 9845.799999999999
$ws.Range("F343").Value = 32
$ws.Range("G343").Value = 2303.04
$ws.Range("F345").Value = 36
$ws.Range("G345").Value = 2210.76
$ws.Range("B346").Value = 24282.41
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 107
$ws.Range("G473").Value = 3512.81
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F486").Value = 73
$ws.Range("G486").Value = 6448.82
$ws.Range("B488").Value = 28809.24
$ws.Range("F508").Value = 52
$ws.Range("G508").Value = 5404.88
$ws.Range("F509").Value = 202
$ws.Range("G509").Value = 16236.76
$ws.Range("B510").Value = 21641.64
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 19
$ws.Range("G573").Value = 776.53
$ws.Range("F578").Value = 64
$ws.Range("G578").Value = 3192.96
$ws.Range("F579").Value = 32
$ws.Range("G579").Value = 2579.2
$ws.Range("F580").Value = 47
$ws.Range("G580").Value = 2678.53
$ws.Range("F582").Value = 23
$ws.Range("G582").Value = 1310.77
$ws.Range("B583").Value = 13466.2
$ws.Range("F599").Value = 1397
$ws.Range("G599").Value = 227864.67
$ws.Range("F601").Value = 374
$ws.Range("G601").Value = 105793.38
$ws.Range("F602").Value = 317
$ws.Range("G602").Value = 45854.05
$ws.Range("B606").Value = 380360.15
$ws.Range("F613").Value = 130
$ws.Range("G613").Value = 20690.8
$ws.Range("B618").Value = 42241
$ws.Range("B619").Value = 1613004.11
$ws.Range("B620").Value = 1613004.11
